$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 38

$ws.Range("E10").Value = 536
$ws.Range("F10").Value = 259
$ws.Range("H10").Value = 353

$ws.Range("E11").Value = 347

$ws.Range("E12").Value = 521
$ws.Range("F12").Value = 281
$ws.Range("H12").Value = 366

$ws.Range("E15").Value = 164

$ws.Range("E16").Value = 202
$ws.Range("F16").Value = 102
$ws.Range("H16").Value = 150

$ws.Range("E18").Value = 52

$ws.Range("E21").Value = 138

$ws.Range("E22").Value = 167
$ws.Range("F22").Value = 89
$ws.Range("H22").Value = 131

$ws.Range("E27").Value = 324
$ws.Range("F27").Value = 162
$ws.Range("H27").Value = 243

$ws.Range("E29").Value = 167

$ws.Range("E37").Value = 157

$ws.Range("E39").Value = 180

$ws.Range("E41").Value = 382
$ws.Range("F41").Value = 179
$ws.Range("H41").Value = 271

$ws.Range("E42").Value = 375
$ws.Range("F42").Value = 207
$ws.Range("H42").Value = 267

$ws.Range("E43").Value = 116
$ws.Range("F43").Value = 62
$ws.Range("H43").Value = 90
